$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF), matching style of existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for column I (I0) and column J (IF), rows 2-13
$dataI = @(6, 2, 1, 1, 1, 1, 1, 1, 1, 1, 3, 1)
$dataJ = @(8, 4, 6, 5, 5, 5, 5, 4, 3, 2, 4, 2)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
